$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44334
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 7000
$ws.Range("P2").Value = 7500
$ws.Range("S2").Value = 2500

# Row 3
$ws.Range("D3").Value = 44334
$ws.Range("M3").Value = 160
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 6500
$ws.Range("S3").Value = 2167

# Row 4
$ws.Range("D4").Value = 44334
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6500
$ws.Range("S4").Value = 2167

# Row 5
$ws.Range("L5").Value = "Tercera"
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 3500
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 3750
$ws.Range("S5").Value = 1250

# Row 6
$ws.Range("D6").Value = 44351
$ws.Range("L6").Value = "Especial"
$ws.Range("N6").Value = 7500
$ws.Range("O6").Value = 8000
$ws.Range("P6").Value = 7750
$ws.Range("S6").Value = 2583

# Row 7
$ws.Range("D7").Value = 44351
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 100
$ws.Range("O7").Value = 6500
$ws.Range("P7").Value = 6250
$ws.Range("S7").Value = 2083

# Row 8
$ws.Range("D8").Value = 44351
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 4500
$ws.Range("O8").Value = 5000
$ws.Range("P8").Value = 4750
$ws.Range("S8").Value = 1583

# Row 9
$ws.Range("D9").Value = 44389
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 7500
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 7750
$ws.Range("S9").Value = 2583

# Row 10
$ws.Range("D10").Value = 44389
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 6500
$ws.Range("S10").Value = 2167

# Row 11
$ws.Range("D11").Value = 44389
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 5500
$ws.Range("O11").Value = 6000
$ws.Range("P11").Value = 5750
$ws.Range("S11").Value = 1917

# Row 24
$ws.Range("D24").Value = 44200
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 4500
$ws.Range("O24").Value = 5000
$ws.Range("P24").Value = 4750
$ws.Range("S24").Value = 1583

# Row 25
$ws.Range("D25").Value = 44200
$ws.Range("M25").Value = 80
$ws.Range("N25").Value = 3500
$ws.Range("O25").Value = 4000
$ws.Range("P25").Value = 3750
$ws.Range("S25").Value = 1250

# Row 26
$ws.Range("D26").Value = 44200
$ws.Range("M26").Value = 120
$ws.Range("N26").Value = 2500
$ws.Range("O26").Value = 3000
$ws.Range("P26").Value = 2750
$ws.Range("S26").Value = 917
